$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue $ws 'D2' '43.147.62'
Set-TextValue $ws 'E2' '  +0.79%  '

# Row 3
Set-TextValue $ws 'D3' '2.368.12'
Set-TextValue $ws 'E3' '  +2.68%  '

# Row 4
Set-TextValue $ws 'D4' '0.999'
Set-TextValue $ws 'E4' '  -0.06%  '

# Row 5
Set-TextValue $ws 'D5' '303.34'
Set-TextValue $ws 'E5' '  +0.54%  '

# Row 6
Set-TextValue $ws 'D6' '96.59'
Set-TextValue $ws 'E6' '  +0.23%  '

# Row 7
Set-TextValue $ws 'D7' '0.509'
Set-TextValue $ws 'E7' '  +0.01%  '

# Row 8
Set-TextValue $ws 'E8' '  -0.06%  '

# Row 9
Set-TextValue $ws 'D9' '0.499'
Set-TextValue $ws 'E9' '  +0.89%  '

# Row 10
Set-TextValue $ws 'D10' '34.19'
Set-TextValue $ws 'E10' '  -1.67%  '

# Row 11
Set-TextValue $ws 'E11' '  +0.05%  '

# Row 12
Set-TextValue $ws 'B12' 'Chainlink'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws 'D12' '18.58'
Set-TextValue $ws 'E12' '  -3.62%  '

# Row 13
Set-TextValue $ws 'B13' 'TRON'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws 'D13' '0.122'
Set-TextValue $ws 'E13' '  +2.69%  '

# Row 14
Set-TextValue $ws 'D14' '6.77'
Set-TextValue $ws 'E14' '  -0.55%  '

# Row 15
Set-TextValue $ws 'D15' '2.732.17'
Set-TextValue $ws 'E15' '  +2.63%  '

# Row 16
Set-TextValue $ws 'D16' '2.361.00'
Set-TextValue $ws 'E16' '  +2.10%  '

# Row 17
Set-TextValue $ws 'D17' '0.800'
Set-TextValue $ws 'E17' '  +1.59%  '

# Row 18
Set-TextValue $ws 'D18' '43.135.98'
Set-TextValue $ws 'E18' '  +0.94%  '

# Row 19
Set-TextValue $ws 'D19' '12.35'
Set-TextValue $ws 'E19' '  +0.63%  '

# Row 20
Set-TextValue $ws 'E20' '  +4.15%  '

# Row 21
Set-TextValue $ws 'D21' '0.0₃0889'
Set-TextValue $ws 'E21' '  -0.48%  '

# Row 22
Set-TextValue $ws 'D22' '68.20'
Set-TextValue $ws 'E22' '  +0.50%  '

# Row 23
Set-TextValue $ws 'D23' '235.98'
Set-TextValue $ws 'E23' '  +0.21%  '

# Row 24
Set-TextValue $ws 'E24' '  -2.88%  '

# Row 25
Set-TextValue $ws 'E25' '  +0.96%  '

# Row 26
Set-TextValue $ws 'E26' '  -0.08%  '

# Row 27
Set-TextValue $ws 'D27' '24.76'
Set-TextValue $ws 'E27' '  +1.19%  '

# Row 28
Set-TextValue $ws 'E28' '  +0.17%  '

# Row 29
Set-TextValue $ws 'D29' '9.15'
Set-TextValue $ws 'E29' '  +0.74%  '

# Row 30
Set-TextValue $ws 'D30' '31.53'
Set-TextValue $ws 'E30' '  -2.34%  '

# Row 31
Set-TextValue $ws 'D31' '0.999'
Set-TextValue $ws 'E31' '  -0.09%  '

# Row 32
Set-TextValue $ws 'E32' '  +1.92%  '

# Row 33
Set-TextValue $ws 'D33' '0.0730'
Set-TextValue $ws 'E33' '  +3.70%  '

# Row 34
Set-TextValue $ws 'D34' '17.27'
Set-TextValue $ws 'E34' '  -2.20%  '

# Row 35
Set-TextValue $ws 'E35' '  +5.06%  '

# Row 36
Set-TextValue $ws 'D36' '4.40'
Set-TextValue $ws 'E36' '  -1.86%  '

# Row 37
Set-TextValue $ws 'E37' '  -1.18%  '

# Row 38
Set-TextValue $ws 'E38' '  +1.27%  '

# Row 39
Set-TextValue $ws 'D39' '2.78'
Set-TextValue $ws 'E39' '  +2.19%  '

# Row 40
Set-TextValue $ws 'D40' '22.47'
Set-TextValue $ws 'E40' '  +10.57%  '

# Row 41
Set-TextValue $ws 'E41' '  -0.14%  '

# Row 42
Set-TextValue $ws 'D42' '1.943.88'
Set-TextValue $ws 'E42' '  -1.65%  '

# Row 43
Set-TextValue $ws 'D43' '102.81'
Set-TextValue $ws 'E43' '  -37.80%  '

# Row 44
Set-TextValue $ws 'E44' '  +0.02%  '

# Row 45
Set-TextValue $ws 'E45' '  +5.12%  '

# Row 46
Set-TextValue $ws 'D46' '9.47'
Set-TextValue $ws 'E46' '  -9.49%  '

# Row 47
Set-TextValue $ws 'D47' '2.75'
Set-TextValue $ws 'E47' '  -1.29%  '

# Row 48
Set-TextValue $ws 'D48' '2.598.16'
Set-TextValue $ws 'E48' '  +2.66%  '

# Row 49
Set-TextValue $ws 'D49' '53.00'
Set-TextValue $ws 'E49' '  -0.88%  '

# Row 50
Set-TextValue $ws 'B50' 'Stacks'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws 'D50' '1.51'
Set-TextValue $ws 'E50' '  +1.24%  '

# Row 51
Set-TextValue $ws 'B51' 'HuobiToken'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D51' '2.78'
Set-TextValue $ws 'E51' '  +0.60%  '
